# Add the new "ODI Bowling Extra" worksheet (scraping additional bowling
# attributes) as the last sheet in the workbook, matching the existing
# "ODI Batting Extra" sheet's layout/style conventions.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Match the outline/page-setup conventions used by every other sheet in
# this workbook (summary rows below detail, summary columns to the right,
# 0.75"/1"/0.5" margins).
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row text
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$col - 1]
}

# Re-use the bold/centered/bordered header style already used by the other
# "Extra" sheet so the new sheet matches the workbook's conventions.
$srcHeader = $wb.Worksheets.Item("ODI Batting Extra").Cells.Item(1, 1)
$srcHeader.Copy()
$ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1,3)).PasteSpecial(-4122)

# MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL for each match row, in the
# same match-code order as the "ODI Batting" / "ODI Batting Extra" sheets.
$rows = @(
    @("3814", "", ""),
    @("3819", "", ""),
    @("3820", "", ""),
    @("3821", "0", "10.00%"),
    @("3822", "", ""),
    @("3836", "0", "60.00%"),
    @("3837", "", ""),
    @("3858", "", ""),
    @("3859", "0", ""),
    @("3863", "0", "10.00%"),
    @("3926", "0", ""),
    @("3928", "0", ""),
    @("4176", "", ""),
    @("4177", "0", ""),
    @("4273", "0", ""),
    @("4274", "1", "10.00%"),
    @("4275", "0", "10.00%"),
    @("4276", "", ""),
    @("4277", "0", ""),
    @("4292", "0", "10.00%")
)

$r = 2
foreach ($row in $rows) {
    for ($col = 1; $col -le 3; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col - 1]
    }
    $r++
}

Write-Output ("Worksheets: " + $wb.Worksheets.Count)
